$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("G5").Value = 1.5
$ws.Range("H5").Value = 3.7
$ws.Range("I5").Value = 6.6
$ws.Range("L5").Value = 1.4
$ws.Range("M5").Value = 2.52
$ws.Range("N5").Value = 2.15
$ws.Range("O5").Value = 1.55
$ws.Range("P5").Value = 1.45
$ws.Range("Q5").Value = 2.4
$ws.Range("R5").Value = 2.32
$ws.Range("S5").Value = 1.47
$ws.Range("T5").Value = 4.9
$ws.Range("U5").Value = 5.7
$ws.Range("V5").Value = 9
$ws.Range("W5").Value = 9.5
$ws.Range("X5").Value = 15
$ws.Range("Y5").Value = 45
$ws.Range("Z5").Value = 7.5
$ws.Range("AA5").Value = 7.6
$ws.Range("AB5").Value = 26
$ws.Range("AC5").Value = 175
$ws.Range("AE5").Value = 13
$ws.Range("AF5").Value = 37
$ws.Range("AG5").Value = 23
$ws.Range("AH5").Value = 175
$ws.Range("AI5").Value = 100
$ws.Range("AJ5").Value = 110

# Row 6
$ws.Range("G6").Value = 2.22
$ws.Range("H6").Value = 3.1
$ws.Range("I6").Value = 3.15
$ws.Range("L6").Value = 1.47
$ws.Range("M6").Value = 2.32
$ws.Range("N6").Value = 2.35
$ws.Range("O6").Value = 1.47
$ws.Range("P6").Value = 1.5
$ws.Range("Q6").Value = 2.25
$ws.Range("R6").Value = 2.05
$ws.Range("S6").Value = 1.6
$ws.Range("T6").Value = 5.9
$ws.Range("U6").Value = 9.25
$ws.Range("V6").Value = 9.75
$ws.Range("W6").Value = 21
$ws.Range("X6").Value = 22
$ws.Range("Y6").Value = 45
$ws.Range("Z6").Value = 6.9
$ws.Range("AA6").Value = 6.2
$ws.Range("AB6").Value = 19
$ws.Range("AC6").Value = 120
$ws.Range("AE6").Value = 7.3
$ws.Range("AF6").Value = 14.5
$ws.Range("AH6").Value = 40
$ws.Range("AJ6").Value = 55

# Row 8
$ws.Range("H8").Value = 3.8
$ws.Range("I8").Value = 4
$ws.Range("L8").Value = 1.24
$ws.Range("M8").Value = 3.3
$ws.Range("N8").Value = 1.72
$ws.Range("O8").Value = 1.9
$ws.Range("R8").Value = 1.7
$ws.Range("S8").Value = 1.91
$ws.Range("T8").Value = 7.6
$ws.Range("U8").Value = 8.5
$ws.Range("Y8").Value = 25
$ws.Range("Z8").Value = 11.75
$ws.Range("AA8").Value = 7.4
$ws.Range("AB8").Value = 15
$ws.Range("AC8").Value = 65
$ws.Range("AE8").Value = 12.5
$ws.Range("AG8").Value = 13.5
$ws.Range("AH8").Value = 60
$ws.Range("AI8").Value = 35
$ws.Range("AJ8").Value = 40

# Row 10
$ws.Range("G10").Value = 1.55
$ws.Range("H10").Value = 3.95
$ws.Range("I10").Value = 5.2
$ws.Range("N10").Value = 1.55
$ws.Range("O10").Value = 2.15
$ws.Range("R10").Value = 1.6
$ws.Range("W10").Value = 12
$ws.Range("X10").Value = 11.5
$ws.Range("Y10").Value = 20
$ws.Range("Z10").Value = 14
$ws.Range("AA10").Value = 8
$ws.Range("AB10").Value = 14
$ws.Range("AC10").Value = 50
$ws.Range("AD10").Value = 300
$ws.Range("AF10").Value = 37

# Row 11
$ws.Range("G11").Value = 4.8
$ws.Range("H11").Value = 4.15
$ws.Range("I11").Value = 1.57
$ws.Range("O11").Value = 2.4
$ws.Range("R11").Value = 1.5
$ws.Range("S11").Value = 2.25
$ws.Range("T11").Value = 19.5
$ws.Range("U11").Value = 35
$ws.Range("V11").Value = 15.5
$ws.Range("W11").Value = 90
$ws.Range("Z11").Value = 17.5
$ws.Range("AA11").Value = 8.75
$ws.Range("AB11").Value = 13
$ws.Range("AC11").Value = 45
$ws.Range("AD11").Value = 250
$ws.Range("AE11").Value = 10.25
$ws.Range("AF11").Value = 9.75
$ws.Range("AH11").Value = 13

# Row 13
$ws.Range("G13").Value = 2.02
$ws.Range("H13").Value = 3.15
$ws.Range("I13").Value = 3.4
$ws.Range("N13").Value = 2
$ws.Range("O13").Value = 1.65
$ws.Range("P13").Value = 1.38
$ws.Range("Q13").Value = 2.47
$ws.Range("T13").Value = 5.8
$ws.Range("U13").Value = 7.9
$ws.Range("V13").Value = 7.3
$ws.Range("W13").Value = 15
$ws.Range("X13").Value = 13.5
$ws.Range("Y13").Value = 23
$ws.Range("Z13").Value = 8.25
$ws.Range("AA13").Value = 5.4
$ws.Range("AB13").Value = 12.5
$ws.Range("AC13").Value = 55
$ws.Range("AD13").Value = 400
$ws.Range("AE13").Value = 7.6
$ws.Range("AF13").Value = 14
$ws.Range("AG13").Value = 10
$ws.Range("AH13").Value = 37
$ws.Range("AI13").Value = 26
$ws.Range("AJ13").Value = 32

# Row 14
$ws.Range("G14").Value = 3.15
$ws.Range("H14").Value = 3.2
$ws.Range("I14").Value = 2.1
$ws.Range("N14").Value = 2.15
$ws.Range("O14").Value = 1.55
$ws.Range("P14").Value = 1.44
$ws.Range("Q14").Value = 2.3
$ws.Range("T14").Value = 6.9
$ws.Range("U14").Value = 12.5
$ws.Range("V14").Value = 9.75
$ws.Range("W14").Value = 32
$ws.Range("X14").Value = 25
$ws.Range("Y14").Value = 35
$ws.Range("Z14").Value = 7.7
$ws.Range("AA14").Value = 5.5
$ws.Range("AB14").Value = 14
$ws.Range("AC14").Value = 70
$ws.Range("AE14").Value = 5.4
$ws.Range("AF14").Value = 7.6
$ws.Range("AG14").Value = 7.7
$ws.Range("AH14").Value = 15
$ws.Range("AI14").Value = 15.5
$ws.Range("AJ14").Value = 28

# Row 15
$ws.Range("G15").Value = 3.4
$ws.Range("H15").Value = 3
$ws.Range("J15").Value = 1.08
$ws.Range("K15").Value = 7.5
$ws.Range("L15").Value = 1.4
$ws.Range("M15").Value = 2.75
$ws.Range("N15").Value = 2.3
$ws.Range("O15").Value = 1.6
$ws.Range("P15").Value = 1.5
$ws.Range("Q15").Value = 2.5
$ws.Range("R15").Value = 1.95
$ws.Range("S15").Value = 1.8
$ws.Range("V15").Value = 13
$ws.Range("Z15").Value = 7.5
$ws.Range("AD15").Value = 351
$ws.Range("AE15").Value = 6.5
$ws.Range("AF15").Value = 10

# Row 17
$ws.Range("G17").Value = 2.25
$ws.Range("H17").Value = 3
$ws.Range("I17").Value = 3.1

# Row 18
$ws.Range("G18").Value = 1.93
$ws.Range("H18").Value = 3.2
$ws.Range("I18").Value = 3.7

# Row 19
$ws.Range("G19").Value = 1.8
$ws.Range("H19").Value = 3.35
$ws.Range("I19").Value = 4.1

# Row 20
$ws.Range("G20").Value = 1.53
$ws.Range("H20").Value = 3.65
$ws.Range("I20").Value = 5.7

# Row 21
$ws.Range("G21").Value = 2
$ws.Range("H21").Value = 3.6
$ws.Range("I21").Value = 3.5
$ws.Range("U21").Value = 9.5
$ws.Range("W21").Value = 17
$ws.Range("AA21").Value = 7
$ws.Range("AF21").Value = 19

# Row 22
$ws.Range("G22").Value = 3.8
$ws.Range("H22").Value = 3.7
$ws.Range("I22").Value = 1.9
$ws.Range("U22").Value = 21
$ws.Range("AH22").Value = 15

# Row 23
$ws.Range("G23").Value = 2.3
$ws.Range("H23").Value = 3.3
$ws.Range("I23").Value = 2.67
$ws.Range("N23").Value = 1.62
$ws.Range("O23").Value = 2.02
$ws.Range("T23").Value = 8.5
$ws.Range("U23").Value = 11
$ws.Range("V23").Value = 7.7
$ws.Range("W23").Value = 20
$ws.Range("X23").Value = 14
$ws.Range("Y23").Value = 18
$ws.Range("Z23").Value = 12.5
$ws.Range("AA23").Value = 5.9
$ws.Range("AB23").Value = 9.75
$ws.Range("AC23").Value = 32
$ws.Range("AD23").Value = 175
$ws.Range("AE23").Value = 9.5
$ws.Range("AF23").Value = 13
$ws.Range("AG23").Value = 8.5
$ws.Range("AH23").Value = 26
$ws.Range("AI23").Value = 16.5
$ws.Range("AJ23").Value = 19

# Row 24
$ws.Range("G24").Value = 1.87
$ws.Range("H24").Value = 3.55
$ws.Range("I24").Value = 3.45
$ws.Range("N24").Value = 1.75
$ws.Range("O24").Value = 1.85
$ws.Range("P24").Value = 1.37
$ws.Range("Q24").Value = 2.5
$ws.Range("T24").Value = 6.6
$ws.Range("U24").Value = 7.8
$ws.Range("V24").Value = 7.2
$ws.Range("W24").Value = 13
$ws.Range("X24").Value = 12
$ws.Range("Y24").Value = 20
$ws.Range("Z24").Value = 11
$ws.Range("AA24").Value = 6.1
$ws.Range("AB24").Value = 12
$ws.Range("AC24").Value = 45
$ws.Range("AD24").Value = 300
$ws.Range("AE24").Value = 9.25
$ws.Range("AF24").Value = 15.5
$ws.Range("AG24").Value = 10
$ws.Range("AH24").Value = 35
$ws.Range("AI24").Value = 24
$ws.Range("AJ24").Value = 28

# Row 28
$ws.Range("G28").Value = 2.55
$ws.Range("I28").Value = 2.8
$ws.Range("W28").Value = 26
$ws.Range("AB28").Value = 17
$ws.Range("AE28").Value = 7
$ws.Range("AG28").Value = 11
$ws.Range("AI28").Value = 26

# Row 32
$ws.Range("J32").Value = 1.05
$ws.Range("K32").Value = 8.5
$ws.Range("N32").Value = 1.93
$ws.Range("O32").Value = 1.88

# Row 36
$ws.Range("H36").Value = 3.6
$ws.Range("I36").Value = 3.25
$ws.Range("N36").Value = 1.75
$ws.Range("O36").Value = 2.05
$ws.Range("T36").Value = 9
$ws.Range("Z36").Value = 13
$ws.Range("AA36").Value = 7
$ws.Range("AF36").Value = 17
$ws.Range("AH36").Value = 34
$ws.Range("AI36").Value = 23
$ws.Range("AJ36").Value = 29
